$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11, column B holds the rule's label (was "R40"); this save renames it
# to the text "1". Go through a text-valued formula + paste-values so the
# literal "1" lands back in the cell as a *string* (t="s") rather than
# being auto-coerced to the number 1, while leaving the cell's existing
# style/format (s="23") untouched.
$cell = $ws.Range("B11")
$cell.Formula = "=""1"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$excel.CutCopyMode = $false
